# ConsumeData.xlsx -> "unify the conception of DataNode, DataTable, Entity"
#
# The sheet that used to describe a generic "Property1" table is renamed to
# "DataNode" (matching the commit message), and the author's cursor/selection
# ends up parked on D42 before the file was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet: "Property1" -> "DataNode"
$ws.Name = "DataNode"

# Move the active selection to D42 (was A9)
$ws.Range("D42").Select()
